# Fixture list update: add missing score entries for rows 50/51 and a new
# fixture row (52) for the Spain vs England match in Berlin, formatting its
# date cell with a short date number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank Home/Away score columns for the last two
# existing fixtures.
$ws.Range("J50").Value = 2
$ws.Range("K50").Value = 1

$ws.Range("J51").Value = 1
$ws.Range("K51").Value = 2

# Append the new fixture as row 52.
$ws.Range("A52").Value = "Sun"

$ws.Range("B52").Value = 45487
$ws.Range("B52").NumberFormat = "d-mmm-yy"

$ws.Range("C52").Value = "21:00:00"
$ws.Range("D52").Value = "Spain"
$ws.Range("G52").Value = "England"
$ws.Range("H52").Value = "Berlin"

# Leave the selection on the newly-added cell, matching where the author's
# cursor ended up.
[void]$ws.Range("H52").Select()
